$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing it to remain text, even when
# the string looks numeric (e.g. "1.00", "0.512") -- mirrors how the
# source data keeps these as literal strings instead of numbers, and
# restores the default "General"/Normal formatting afterwards so no
# stray number-format/style is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '61.075.31'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '2.969.53'
$ws.Range('E3').Value = '  -1.19%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  +0.09%  '
Set-TextValue $ws.Range('D5') '592.91'
$ws.Range('E5').Value = '  +1.43%  '
Set-TextValue $ws.Range('D6') '141.31'
$ws.Range('E6').Value = '  -3.63%  '
$ws.Range('E7').Value = '  +0.13%  '
Set-TextValue $ws.Range('D8') '0.512'
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').Value = '2.964.78'
$ws.Range('E9').Value = '  -1.32%  '
Set-TextValue $ws.Range('D10') '0.144'
$ws.Range('E10').Value = '  -2.99%  '
$ws.Range('E11').Value = '  +3.89%  '
Set-TextValue $ws.Range('D12') '0.449'
$ws.Range('E12').Value = '  +1.59%  '
Set-TextValue $ws.Range('D13') '0.0000225'
$ws.Range('E13').Value = '  -1.52%  '
Set-TextValue $ws.Range('D14') '33.90'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').Value = '3.464.23'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').Value = '61.223.89'
$ws.Range('E17').Value = '  -1.91%  '
Set-TextValue $ws.Range('D18') '6.81'
$ws.Range('E18').Value = '  -2.96%  '
$ws.Range('D19').Value = '2.977.87'
$ws.Range('E19').Value = '  -0.91%  '
Set-TextValue $ws.Range('D20') '447.18'
$ws.Range('E20').Value = '  -2.58%  '
Set-TextValue $ws.Range('D21') '14.04'
$ws.Range('E21').Value = '  +1.04%  '
Set-TextValue $ws.Range('D22') '0.676'
$ws.Range('E22').Value = '  -0.53%  '
Set-TextValue $ws.Range('D23') '7.21'
$ws.Range('E23').Value = '  -2.14%  '
Set-TextValue $ws.Range('D24') '82.04'
$ws.Range('E24').Value = '  +2.44%  '
Set-TextValue $ws.Range('D25') '2.14'
$ws.Range('E25').Value = '  -5.62%  '
Set-TextValue $ws.Range('D26') '11.88'
$ws.Range('E26').Value = '  -2.76%  '
Set-TextValue $ws.Range('D27') '10.23'
$ws.Range('E27').Value = '  +1.37%  '
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D29') '1.00'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D30') '2.65'
$ws.Range('E30').Value = '  +1.43%  '
Set-TextValue $ws.Range('D31') '6.99'
$ws.Range('E31').Value = '  -3.00%  '
Set-TextValue $ws.Range('D32') '2.02'
$ws.Range('E32').Value = '  -3.07%  '
Set-TextValue $ws.Range('D33') '27.05'
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('D35').Value = '0.0₃0797'
$ws.Range('E35').Value = '  +0.99%  '
Set-TextValue $ws.Range('D36') '1.00'
$ws.Range('E36').Value = '  -2.45%  '
Set-TextValue $ws.Range('D37') '5.72'
$ws.Range('E37').Value = '  -0.30%  '
Set-TextValue $ws.Range('D38') '50.07'
$ws.Range('E38').Value = '  +0.03%  '
Set-TextValue $ws.Range('D39') '2.03'
$ws.Range('E39').Value = '  -4.02%  '
Set-TextValue $ws.Range('D40') '8.93'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('E41').Value = '  +7.60%  '
Set-TextValue $ws.Range('D42') '2.81'
$ws.Range('E42').Value = '  -4.70%  '
Set-TextValue $ws.Range('D43') '386.87'
$ws.Range('E43').Value = '  -5.57%  '
Set-TextValue $ws.Range('D44') '0.0346'
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Range('D45') '38.23'
$ws.Range('E45').Value = '  -2.18%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D46') '0.264'
$ws.Range('E46').Value = '  -4.54%  '
$ws.Range('D47').Value = '2.682.52'
$ws.Range('E47').Value = '  -3.20%  '
Set-TextValue $ws.Range('D48') '129.04'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('E49').Value = '  +0.14%  '
Set-TextValue $ws.Range('D50') '0.106'
$ws.Range('E50').Value = '  -1.47%  '
Set-TextValue $ws.Range('D51') '2.11'
$ws.Range('E51').Value = '  -1.23%  '
